# Update the "Ngày cập nhật" (updated_at) timestamp column on the
# "Products" sheet with a fresh backup timestamp (2025-09-23 run #3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

$ws.Range("U3").Value2 = "2025-09-23T13:37:31.130632"
$ws.Range("U4").Value2 = "2025-09-23T13:37:31.132284"
$ws.Range("U5").Value2 = "2025-09-23T13:37:31.133203"
$ws.Range("U6").Value2 = "2025-09-23T13:37:31.134036"
$ws.Range("U7").Value2 = "2025-09-23T13:37:31.134926"
